$d = $word.ActiveDocument

$found1 = $d.Content.Find.Execute('Ativação: 01/01/1996', $true, $false, $false, $false, $false, $true, 1, $false, 'Ativação: 01/01/2022', 2)
Write-Host "Step1 (Ativacao):" $found1

# Step 2: split Docente paragraph into two runs
$pDocente = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13,[char]7) -eq "519033 - Carlos Yujiro Shigue") {
        $pDocente = $d.Paragraphs.Item($i)
        break
    }
}
if ($pDocente -eq $null) { Write-Host "Docente paragraph not found!" } else {
    $rDoc = $pDocente.Range
    $rDoc2 = $d.Range($rDoc.Start, $rDoc.End - 1)
    $rDoc2.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>5840897 - Clodoaldo Saron</w:t><w:br/></w:r><w:r><w:t>1033242 - Fábio Herbst Florenzano</w:t></w:r></w:p></w:body></w:document>')
    Write-Host "Step2 (Docente split) done"
}

$found3 = $d.Content.Find.Execute('A avaliação será feita por meio de provas escritas.', $true, $false, $false, $false, $false, $true, 1, $false, 'A avaliação será feita por meio de Provas Escritas, Estudos de Casos e Desenvolvimento de Projetos, sendo necessário utilizar pelo menos dois critérios de avaliação diferentes.', 2)
Write-Host "Step3 (Metodo):" $found3

$found4 = $d.Content.Find.Execute('A Nota final (NF) será calculada da seguinte maneira:NF = (P1 + 2*P2)/3', $true, $false, $false, $false, $false, $true, 1, $false, 'A Nota final (NF) será calculada da seguinte maneira: NF = (P+EC+Projetos)/3', 2)
Write-Host "Step4 (Criterio):" $found4

$found5 = $d.Content.Find.Execute('A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2', $true, $false, $false, $false, $false, $true, 1, $false, 'Não consta recuperação', 2)
Write-Host "Step5 (Norma recuperacao):" $found5

# Step 6: Bibliografia - locate paragraph starting with "1. J. Margolis."
$pBib = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("1. J. Margolis.")) {
        $pBib = $d.Paragraphs.Item($i)
        break
    }
}
if ($pBib -eq $null) { Write-Host "Bibliografia paragraph not found!" } else {
    $rBib = $pBib.Range
    $rBib2 = $d.Range($rBib.Start, $rBib.End - 1)
    $rBib2.Text = '1. J. Margolis. Engineering Plastics Handbook. McGraw-Hill Professional, 2005. 2. Nigel Mills. Plastics - Microstructure and Engineering Applications. Butterworth-Heineman, 2005. 3. Walter Michaeli, TEcnologia dos Plasticos. Ed. Blucher 4. Hélio Wiebeck, Júlio Harada. Plásticos de Engenharia - Tecnologia e Aplicações. São Paulo: Editora Artliber, 2005. 5. E. B. Mano, L. C. Mendes. Identificação de Plásticos, Borrachas e Fibras. São Paulo: Editora Edgard Blucher, 2000. 6. Marcelo Rabello. Aditivação de Polímeros. São Paulo: Editora Artliber, 2004. 7. Jan C.J. Bart. Additives in Polymers. New York: John Wiley & Sons, 2005. 8. Marino Xanthos. Functional Fillers for Plastics. Wiley-VCH Verlag GmbH, 2005. 9. Silvio Manrich. Processamento de Termoplásticos. Editora Artliber, 2005. 10. G.H. Michler, F.J. Baltá-Calleja. Mechanical Properties of Polymers Based on Nanostructure and Morphology. Boca Raton: CRC Press, 2005. 11. A. M. Piva, H. Wiebeck. Reciclagem do P. São Paulo: Editora Artliber". Manas Chanda, ,Salil K. Roy  Plastics Fabrication and Recycling'
    Write-Host "Step6 (Bibliografia) done"
}

